$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new new sheet")

$ws.Range("C2").Value = 94547
$ws.Range("C3").Value = 59346
$ws.Range("C4").Value = 83076
$ws.Range("C5").Value = 26068
$ws.Range("C6").Value = 55861
$ws.Range("C7").Value = 28547
$ws.Range("C8").Value = 11840
$ws.Range("C9").Value = 41102
$ws.Range("C10").Value = 87681
$ws.Range("C11").Value = 15646
$ws.Range("C12").Value = 35703
$ws.Range("C13").Value = 87336
